# Update Leave Card 12/22/2023 10:59 AM
# Applies the monthly leave-credit rollover: shifts the SL "period" dates in
# column A down to end-of-month values continuing the series, and fills in
# the next six months of earned Sick Leave credits (1.25 each) so the
# calculated BALANCE columns recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")

# --- Column A: PERIOD dates (end-of-month series, continuing from row 57) ---
$ws.Range("A57").Value = "01/31/2023"
$ws.Range("A58").Value = "02/28/2023"
$ws.Range("A59").Value = "03/31/2023"
$ws.Range("A60").Value = "04/30/2023"
$ws.Range("A61").Value = "05/31/2023"
$ws.Range("A62").Value = "06/30/2023"
$ws.Range("A63").Value = "07/31/2023"
$ws.Range("A64").Value = "08/31/2023"
$ws.Range("A65").Value = "09/30/2023"
$ws.Range("A66").Value = "10/31/2023"
$ws.Range("A67").Value = "11/30/2023"
$ws.Range("A68").Value = "12/31/2023"
$ws.Range("A69").Value = "01/31/2024"
$ws.Range("A70").Value = "02/29/2024"
$ws.Range("A71").Value = "03/31/2024"
$ws.Range("A72").Value = "04/30/2024"
$ws.Range("A73").Value = "05/31/2024"
$ws.Range("A74").Value = "06/30/2024"
$ws.Range("A75").Value = "07/31/2024"
$ws.Range("A76").Value = "08/31/2024"
$ws.Range("A77").Value = "09/30/2024"
$ws.Range("A78").Value = "10/31/2024"
$ws.Range("A79").Value = "11/30/2024"
$ws.Range("A80").Value = "12/31/2024"
$ws.Range("A81").Value = "01/31/2025"
$ws.Range("A82").Value = "02/28/2025"
$ws.Range("A83").Value = "03/31/2025"
$ws.Range("A84").Value = "04/30/2025"
$ws.Range("A85").Value = "05/31/2025"
$ws.Range("A86").Value = "06/30/2025"
$ws.Range("A87").Value = "07/31/2025"
$ws.Range("A88").Value = "08/31/2025"
$ws.Range("A89").Value = "09/30/2025"
$ws.Range("A90").Value = "10/31/2025"
$ws.Range("A91").Value = "11/30/2025"

# --- Column C: EARNED (Sick Leave) credits for the newly-rolled months ---
$ws.Range("C60").Value = 1.25
$ws.Range("C61").Value = 1.25
$ws.Range("C62").Value = 1.25
$ws.Range("C63").Value = 1.25
$ws.Range("C64").Value = 1.25
$ws.Range("C65").Value = 1.25

$excel.CalculateFullRebuild()
